$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 12 data for Tinder / MTCH (adding MTCH to the portfolio)
$ws.Range("A12").Value = "Tinder"
$ws.Range("B12").Value = "MTCH"
$ws.Range("C12").Value = "Market Returns"

# Update the date range string for all data rows (2-12) to the refreshed window
$ws.Range("D2:D12").Value = "2022-10-03 - 2024-09-20"

# Update Beta / Standard Error / t-value / p-value for each row (refreshed Sp500 betas, plus new MTCH row)
$ws.Range("E2").Value = 1.004615391668081
$ws.Range("F2").Value = 0.00995854079907098
$ws.Range("G2").Value = 100.8797786681559
$ws.Range("H2").Value = [double]"1.07173329996564E-29"
$ws.Range("E3").Value = 1.310756136243378
$ws.Range("F3").Value = 0.3588983204424652
$ws.Range("G3").Value = 3.652165701492895
$ws.Range("H3").Value = 0.001487785739533245
$ws.Range("E4").Value = 1.445110674142892
$ws.Range("F4").Value = 0.5814516207294718
$ws.Range("G4").Value = 2.485349808346736
$ws.Range("H4").Value = 0.02144572449122563
$ws.Range("E5").Value = 2.150786540329463
$ws.Range("F5").Value = 0.4930045616409439
$ws.Range("G5").Value = 4.362609816774647
$ws.Range("H5").Value = 0.0002730074041201833
$ws.Range("E6").Value = 1.244794109883067
$ws.Range("F6").Value = 0.4502407874277962
$ws.Range("G6").Value = 2.764729772694553
$ws.Range("H6").Value = 0.01161085469338004
$ws.Range("E7").Value = 1.808481337309522
$ws.Range("F7").Value = 0.4536861211055292
$ws.Range("G7").Value = 3.986194977493839
$ws.Range("H7").Value = 0.0006714998236065933
$ws.Range("E8").Value = 1.820325349745084
$ws.Range("F8").Value = 0.8016265339072233
$ws.Range("G8").Value = 2.270789791441409
$ws.Range("H8").Value = 0.03380331770559451
$ws.Range("E9").Value = 1.394180262023718
$ws.Range("F9").Value = 0.3565260215679585
$ws.Range("G9").Value = 3.910458641678611
$ws.Range("H9").Value = 0.0008045976145165426
$ws.Range("E10").Value = 1.110218960991679
$ws.Range("F10").Value = 0.2922566511359825
$ws.Range("G10").Value = 3.798780820475191
$ws.Range("H10").Value = 0.00105002869816451
$ws.Range("E11").Value = 0.8432122289204602
$ws.Range("F11").Value = 0.4147805513157619
$ws.Range("G11").Value = 2.032911683649661
$ws.Range("H11").Value = 0.05489635629588587
$ws.Range("E12").Value = 2.292724702980512
$ws.Range("F12").Value = 0.5988664939774522
$ws.Range("G12").Value = 3.828440438791413
$ws.Range("H12").Value = 0.0009784048750956023
